# Updated symbol list on Sat Dec 24 11:19:35 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking quotes as plain text
# (t="inlineStr" in the source workbook), so each update forces the cell
# to Text (via NumberFormat "@") before writing the new quote, then
# restores the original "Normal" style so no stray number-format style
# sticks around on the cell. Column E text-only updates are plain
# string writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue "D2"  "244.76"
Set-TextValue "D3"  "21.96"
Set-TextValue "D4"  "5.400"
Set-TextValue "D5"  "0.06008"
Set-TextValue "D7"  "0.8105"
Set-TextValue "D8"  "0.9542"
Set-TextValue "D9"  "0.1425"
Set-TextValue "D10" "0.07395"
Set-TextValue "D11" "0.03376"
Set-TextValue "D12" "0.03053"
Set-TextValue "D13" "0.09421"
Set-TextValue "D14" "4.003"
Set-TextValue "D15" "0.001587"
Set-TextValue "D16" "0.04824"
Set-TextValue "D18" "0.006133"
Set-TextValue "D19" "0.005046"
Set-TextValue "D20" "0.0009906"
Set-TextValue "D22" "3.696"
Set-TextValue "D23" "6.404"
Set-TextValue "D24" "2.186"
Set-TextValue "D41" "0.006572"
Set-TextValue "D42" "0.1074"
Set-TextValue "D43" "0.002901"
Set-TextValue "D44" "0.005233"
Set-TextValue "D45" "0.00005214"
Set-TextValue "D48" "0.02022"

# Column E (Volume(1h) label) updates
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
